$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 419.0909
$ws.Range("I15").Value = 419.0909
$ws.Range("K15").Value = 1257.2727
$ws.Range("M15").Value = -1088.2727
$ws.Range("H19").Value = 2992.75
$ws.Range("I19").Value = 2611.25
$ws.Range("K19").Value = 2611.25
$ws.Range("M19").Value = -2436.25
$ws.Range("H62").Value = 3629
$ws.Range("I62").Value = 2999.5
$ws.Range("K62").Value = 2999.5
$ws.Range("M62").Value = -2375.5
$ws.Range("H65").Value = 3629
$ws.Range("I65").Value = 2999.5
$ws.Range("K65").Value = 14997.5
$ws.Range("M65").Value = -11877.5
$ws.Range("H70").Value = 2507.9167
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 2507.9167
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H86").Value = 5565.533
$ws.Range("J86").Value = 5837.625
$ws.Range("L86").Value = 5837.625
$ws.Range("N86").Value = -8083.625
$ws.Range("H89").Value = 5565.533
$ws.Range("J89").Value = 5837.625
$ws.Range("L89").Value = 29188.125
$ws.Range("N89").Value = -40420.125
$ws.Range("H137").Value = 1335.8695
$ws.Range("I137").Value = 1398
$ws.Range("J137").Value = 1278.9166
$ws.Range("K137").Value = 4194
$ws.Range("L137").Value = 3836.7498
$ws.Range("M137").Value = -1644
$ws.Range("N137").Value = -8936.7498
$ws.Range("H138").Value = 3488.5715
$ws.Range("I138").Value = 3148.889
$ws.Range("J138").Value = 4100
$ws.Range("K138").Value = 9446.667000000001
$ws.Range("L138").Value = 12300
$ws.Range("M138").Value = -4306.667000000001
$ws.Range("N138").Value = -22580
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 807.5
$ws.Range("I74").Value = 807.5
$ws.Range("K74").Value = 807.5
$ws.Range("M74").Value = 66.5
$ws.Range("H77").Value = 807.5
$ws.Range("I77").Value = 807.5
$ws.Range("K77").Value = 4037.5
$ws.Range("M77").Value = 330.5
$ws.Range("H110").Value = 1270
$ws.Range("I110").Value = 1211
$ws.Range("K110").Value = 1211
$ws.Range("M110").Value = 834

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2468.7273
$ws.Range("I99").Value = 2329.5715
$ws.Range("K99").Value = 2329.5715
$ws.Range("M99").Value = -831.5715
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2990
$ws.Range("J15").Value = 2990
$ws.Range("L15").Value = 2990
$ws.Range("N15").Value = -3330
$ws.Range("H31").Value = 2612.25
$ws.Range("I31").Value = 1350
$ws.Range("K31").Value = 1350
$ws.Range("M31").Value = -1055
$ws.Range("H34").Value = 2612.25
$ws.Range("I34").Value = 1350
$ws.Range("K34").Value = 1350
$ws.Range("M34").Value = -1148
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H105").Value = 13727.5
$ws.Range("I105").Value = 20886
$ws.Range("J105").Value = 1796.6666
$ws.Range("K105").Value = 20886
$ws.Range("L105").Value = 1796.6666
$ws.Range("M105").Value = -19139
$ws.Range("N105").Value = -5290.6666
$ws.Range("H134").Value = 2799.5
$ws.Range("I134").Value = 2066
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6198
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3663
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("H141").Value = 100024.695
$ws.Range("I141").Value = 69999.5
$ws.Range("J141").Value = 105483.82
$ws.Range("K141").Value = 69999.5
$ws.Range("L141").Value = 105483.82
$ws.Range("M141").Value = -64819.5
$ws.Range("N141").Value = -115843.82
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 525.5
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 534
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1602
$ws.Range("M12").Value = -1327
$ws.Range("N12").Value = -1948
$ws.Range("H41").Value = 450
$ws.Range("J41").Value = 600
$ws.Range("L41").Value = 1800
$ws.Range("N41").Value = -2476
$ws.Range("H60").Value = 7300
$ws.Range("I60").Value = 6950
$ws.Range("K60").Value = 20850
$ws.Range("M60").Value = -20599
$ws.Range("H61").Value = 443.33334
$ws.Range("I61").Value = 200
$ws.Range("J61").Value = 565
$ws.Range("K61").Value = 600
$ws.Range("L61").Value = 1695
$ws.Range("M61").Value = -385
$ws.Range("N61").Value = -2125
$ws.Range("H80").Value = 7426.5713
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 8164.3335
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 24493.0005
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -26365.0005
$ws.Range("H83").Value = 7426.5713
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 8164.3335
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 73479.0015
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -82839.0015
$ws.Range("H137").Value = 3415.3914
$ws.Range("I137").Value = 1567.5
$ws.Range("J137").Value = 3804.4211
$ws.Range("K137").Value = 4702.5
$ws.Range("L137").Value = 11413.2633
$ws.Range("M137").Value = 397.5
$ws.Range("N137").Value = -21613.2633
$ws.Range("H138").Value = 2498.3333
$ws.Range("I138").Value = 2247.5
$ws.Range("K138").Value = 6742.5
$ws.Range("M138").Value = -1602.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 2400
$ws.Range("I5").Value = 2400
$ws.Range("K5").Value = 2400
$ws.Range("M5").Value = -2287
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15686
$ws.Range("H93").Value = 848.625
$ws.Range("I93").Value = 798.1667
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 798.1667
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 449.8333
$ws.Range("N93").Value = -3496
$ws.Range("H132").Value = 3258.6667
$ws.Range("I132").Value = 3258.6667
$ws.Range("K132").Value = 9776.000100000001
$ws.Range("M132").Value = -7246.000100000001
$ws.Range("H136").Value = 2127.3
$ws.Range("I136").Value = 2216.5
$ws.Range("J136").Value = 1770.5
$ws.Range("K136").Value = 6649.5
$ws.Range("L136").Value = 5311.5
$ws.Range("M136").Value = -4099.5
$ws.Range("N136").Value = -10411.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H96").Value = 1813.2667
$ws.Range("J96").Value = 1680
$ws.Range("L96").Value = 1680
$ws.Range("N96").Value = -4426
$ws.Range("H132").Value = 1794.2354
$ws.Range("I132").Value = 1579.0714
$ws.Range("K132").Value = 4737.2142
$ws.Range("M132").Value = -2207.2142
$ws.Range("H136").Value = 691
$ws.Range("I136").Value = 691
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2073
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 477
$ws.Range("M29").ClearContents()
$ws.Range("N136").ClearContents()
